# Reduce adobe (MUR+ADO) and increase confined masonry (MCF) shares in the
# "Trade" (commercial) building-distribution cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Trade" column's multi-line distribution text lives in C2.
$tradeText = $ws.Range("C2").Value2

$tradeText = $tradeText.Replace("18.0% MUR+ADO/LWAL+CDN/H:1/Trade", "8.0% MUR+ADO/LWAL+CDN/H:1/Trade")
$tradeText = $tradeText.Replace("20.0% MCF/LWAL+CDL/H:1/Trade", "25.0% MCF/LWAL+CDL/H:1/Trade")
$tradeText = $tradeText.Replace("20.0% MCF/LWAL+CDN/H:1/Trade", "25.0% MCF/LWAL+CDN/H:1/Trade")

$ws.Range("C2").Value2 = $tradeText

# Restore the selection recorded after the edit.
$ws.Range("B7").Select()
